$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Chiro Survival" sheet right before "Baetid Survival Rates"
# ---------------------------------------------------------------------------
$baetidSurvival = $wb.Worksheets.Item("Baetid Survival Rates")
$chiro = $wb.Worksheets.Add($baetidSurvival)
$chiro.Name = "Chiro Survival"

$chiro.Range("A1").Value = "Temp"
$chiro.Range("B1").Value = "Survival"
$chiro.Range("C1").Value = "Citation"

$chiroData = @(
    @(5, 16.43762508),
    @(5, 0),
    @(9, 44.807415470000002),
    @(9, 0.30897791499999999),
    @(14, 69.146448509999999),
    @(14, 84.645904290000004),
    @(18, 54.50651311),
    @(18, 9.5137108949999991)
)

$r = 2
foreach ($row in $chiroData) {
    $chiro.Cells.Item($r, 1).Value = $row[0]
    $chiro.Cells.Item($r, 2).Value = $row[1]
    $chiro.Cells.Item($r, 3).Value = "Eggermont & Heiri, 2012"
    $r++
}

# ---------------------------------------------------------------------------
# 2. Append a new "Sheet4" at the very end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)
$sheet4.Name = "Sheet4"

$sheet4.Range("A1").Value = "Taxa"
$sheet4.Range("B1").Value = "Max Discharge m3/s"
$sheet4.Range("C1").Value = "Mortality"
$sheet4.Range("D1").Value = "Citation"
$sheet4.Range("F1").Value = "Bankful Discharge"
$sheet4.Range("G1").Value = "Max/Bankful"

# B, C(mortality), D(citation), F(bankful discharge)
$chironRows = @(
    @(0.3,                    0.84099999999999997, "Boulton et al., 1992", 36.8119011, $null),
    @(1.7,                    0.98099999999999998, "Boulton et al., 1992", 36.8119011, $null),
    @(8.5,                    0.95599999999999996, "Boulton et al., 1992", 36.8119011, $null),
    @(26,                     0.90100000000000002, "Boulton et al., 1992", 36.8119011, $null),
    @(0.54400000000000004,    0.747,                "Cobb et al, 1992",     2.7,        $null),
    @(0.54400000000000004,    0,                    "Cobb et al, 1992",     2.7,        "(no mortality - increase)"),
    @(0.54400000000000004,    0.49,                 "Cobb et al, 1992",     1.6,        $null),
    @(0.47099999999999997,    0.89300000000000002, "Cobb et al, 1992",     2.7,        $null),
    @(0.47099999999999997,    0.49,                 "Cobb et al, 1992",     2.7,        $null),
    @(0.47099999999999997,    0,                    "Cobb et al, 1992",     1.6,        "(no mortality - increase)")
)

$r = 2
foreach ($row in $chironRows) {
    $sheet4.Cells.Item($r, 1).Value = "Chironomidae"
    $sheet4.Cells.Item($r, 2).Value = $row[0]
    $sheet4.Cells.Item($r, 3).Value = $row[1]
    $sheet4.Cells.Item($r, 3).NumberFormat = "0%"
    $sheet4.Cells.Item($r, 4).Value = $row[2]
    $sheet4.Cells.Item($r, 6).Value = $row[3]
    $sheet4.Cells.Item($r, 7).Formula = "=B$r/F$r"
    if ($row[4]) {
        $sheet4.Cells.Item($r, 8).Value = $row[4]
    }
    $r++
}
